# Adds a new 4th data row ("Le Van C") to the Sheet1 user-import example,
# mirroring rows 2 and 3: plain names, SSN-like numeric-text fields, a
# birthday date, an email hyperlink, an image-url hyperlink repeated in
# three columns and the trailing boolean-ish flag columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imgUrl = "https://images.pexels.com/photos/2071882/pexels-photo-2071882.jpeg?auto=compress&cs=tinysrgb&dpr=1&w=500"

# Values are written left-to-right (column order) so that any brand new
# shared-string entries get appended to xl/sharedStrings.xml in the same
# order the source workbook used.

# A4 / B4 / C4 - plain text (A/B reuse existing shared strings)
$ws.Range("A4").Value = "Lê"
$ws.Range("B4").Value = "Văn"
$ws.Range("C4").Value = "C"

# D4 / E4 - numeric-looking identifiers that must stay text, like row 3
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "123456789005"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "123456005"

# F4 - birthday, stored as a real date serial number (same as rows 2-3)
$ws.Range("F4").Value = 36526

# G4 - plain text
$ws.Range("G4").Value = "test"

# H4 - e-mail address with a mailto hyperlink
$ws.Range("H4").Value = "test456@gmail.com"

# I4 - real number (password/pin), unlike the text columns around it
$ws.Range("I4").Value = 123456

# J4 - phone number, kept as text so the leading zero survives
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "0123456789"

# K4 / L4 / M4 - same image URL repeated, each with its own hyperlink
$ws.Range("K4").Value = $imgUrl
$ws.Range("L4").Value = $imgUrl
$ws.Range("M4").Value = $imgUrl

# N4 / P4 - text flags "0"/"1"; O4 - real numeric 0
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = "0"

$ws.Range("O4").Value = 0

$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "1"

# Wire up the hyperlinks for the e-mail and the three image columns.
$ws.Hyperlinks.Add($ws.Range("H4"), "mailto:test456@gmail.com")
$ws.Hyperlinks.Add($ws.Range("K4"), $imgUrl)
$ws.Hyperlinks.Add($ws.Range("L4"), $imgUrl)
$ws.Hyperlinks.Add($ws.Range("M4"), $imgUrl)

# Re-apply row 3's cell formatting onto row 4 so the new row reuses the
# same styles as the existing data rows instead of creating new ones.
$ws.Range("A3:P3").Copy()
$ws.Range("A4:P4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
